$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.466.68"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.625.74"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'212.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.498"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.0622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'18.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.852.27"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.632.44"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'64.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "26.475.92"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'213.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'4.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").Value = "'148.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "1.219.72"
$ws.Range("E35").Value = "  +4.25%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "'0.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").Value = "'2.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.11%  "
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "'5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "1.762.66"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "'92.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'1.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "'54.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "'0.0510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.407"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  +0.39%  "
